# Update EffectiveDate values for the two data rows (F2 and F3)
# from "07302023" to "08152023", and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = "08152023"
$ws.Range("F3").Value = "08152023"

# Update the saved selection/active cell to F8 (matches the target sheetView)
$ws.Range("F8").Select()
